$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    "datastruct",
    "algo",
    "sysprog",
    "ver",
    "build",
    "test",
    "probdec",
    "sysdec",
    "com",
    "orgfile",
    "ordxfile",
    "tree",
    "read",
    "def",
    "err",
    "ide",
    "api",
    "fw",
    "req",
    "scr",
    "db",
    "lang",
    "plat",
    "yrs",
    "dom",
    "tool",
    "langexp",
    "cbaseexp",
    "upcom",
    "platint",
    "book",
    "blog"
)

$row = 2
foreach ($val in $values) {
    $ws.Range("B$row").Value = $val
    $row = $row + 1
}

$ws.Activate()
$ws.Range("B15").Select()
